# Updates Price (D) and Volume(1h) (E) columns with refreshed crypto-market
# data, matching the source GitHub Actions scrape commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.936.44'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '2.337.48'
$ws.Range('E3').Value = '  +0.75%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = "'302.31"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').Value = "'94.25"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.37%  '
$ws.Range('D7').Value = "'0.501"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.94%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = "'0.494"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.24%  '
$ws.Range('D10').Value = "'33.93"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.96%  '
$ws.Range('D11').Value = "'0.0782"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.37%  '
$ws.Range('D12').Value = "'18.58"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.53%  '
$ws.Range('E13').Value = '  +1.42%  '
$ws.Range('D14').Value = "'6.72"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.19%  '
$ws.Range('D15').Value = '2.705.69'
$ws.Range('E15').Value = '  +0.85%  '
$ws.Range('D16').Value = '2.324.50'
$ws.Range('E16').Value = '  -0.83%  '
$ws.Range('D17').Value = "'0.794"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.78%  '
$ws.Range('D18').Value = '42.887.37'
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('D19').Value = "'12.06"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.95%  '
$ws.Range('D20').Value = "'6.18"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.30%  '
$ws.Range('D21').Value = '0.0₃0887'
$ws.Range('E21').Value = '  -0.95%  '
$ws.Range('D22').Value = "'67.88"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').Value = "'235.24"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.61%  '
$ws.Range('D24').Value = "'2.21"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.17%  '
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('E26').Value = '  -1.48%  '
$ws.Range('D27').Value = "'24.53"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.56%  '
$ws.Range('E28').Value = '  +13.81%  '
$ws.Range('D29').Value = "'9.16"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.70%  '
$ws.Range('D30').Value = "'31.38"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.41%  '
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('D32').Value = "'4.98"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.11%  '
$ws.Range('D33').Value = "'0.0735"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.08%  '
$ws.Range('D34').Value = "'17.19"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.80%  '
$ws.Range('D35').Value = "'1.82"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.15%  '
$ws.Range('D36').Value = "'4.36"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.62%  '
$ws.Range('E37').Value = '  -1.82%  '
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('D39').Value = "'121.55"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -26.05%  '
$ws.Range('D40').Value = "'2.75"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.06%  '
$ws.Range('D41').Value = "'22.10"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +16.17%  '
$ws.Range('E42').Value = '  -1.13%  '
$ws.Range('D43').Value = '1.934.94'
$ws.Range('E43').Value = '  -2.34%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').Value = "'10.06"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.69%  '
$ws.Range('D46').Value = "'2.10"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.07%  '
$ws.Range('D47').Value = "'2.70"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.88%  '
$ws.Range('D48').Value = '2.570.09'
$ws.Range('E48').Value = '  +0.85%  '
$ws.Range('E49').Value = '  +0.85%  '
$ws.Range('D50').Value = "'52.73"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.77%  '
$ws.Range('D51').Value = "'71.62"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.96%  '
